# LogInTest data update
# Rewrites the LogIn sheet test data: username/password pairs for
# correct and invalid login scenarios.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Correct Username"
$ws.Range("B1").Value = "admin"

$ws.Range("A2").Value = "Correct Password"
$ws.Range("B2").Value = "admin"

$ws.Range("A4").Value = "Invaid Username"
$ws.Range("B4").Value = "invalidUsername"

$ws.Range("A5").Value = "Invalid Password"
$ws.Range("B5").Value = "invalidPassword"

$ws.Range("B11").Select()
